$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: build the subscript-three character (U+2083) safely so PowerShell
# does not coerce a [char] + numeric-looking-string "+" into an integer add.
$sub3 = [string][char]0x2083

# Price (column D) values are plain text in the workbook (some contain two
# dots as thousands separators, e.g. "51.293.84", so they are never valid
# numbers -- but plain decimals like "0.999" or "1.00" WOULD be auto-converted
# to numeric values by Excel, losing the original text formatting/precision).
# To preserve them as literal text we temporarily force a Text number format,
# write the value, then restore the cell to the default "Normal" style so no
# stray formatting is left behind.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "51.293.84"
$ws.Range("E2").Value = "  +0.51%  "

Set-TextValue $ws.Range("D3") "2.979.41"
$ws.Range("E3").Value = "  +1.57%  "

Set-TextValue $ws.Range("D4") "0.999"
$ws.Range("E4").Value = "  -0.03%  "

Set-TextValue $ws.Range("D5") "381.09"
$ws.Range("E5").Value = "  +2.02%  "

Set-TextValue $ws.Range("D6") "103.06"
$ws.Range("E6").Value = "  +2.49%  "

Set-TextValue $ws.Range("D7") "0.546"
$ws.Range("E7").Value = "  +2.37%  "

Set-TextValue $ws.Range("D9") "0.592"
$ws.Range("E9").Value = "  +1.79%  "

Set-TextValue $ws.Range("D10") "36.67"
$ws.Range("E10").Value = "  +1.61%  "

$ws.Range("E11").Value = "  -0.85%  "

Set-TextValue $ws.Range("D12") "0.0860"

Set-TextValue $ws.Range("D13") "3.447.36"
$ws.Range("E13").Value = "  +1.23%  "

Set-TextValue $ws.Range("D14") "7.82"
$ws.Range("E14").Value = "  +4.68%  "

Set-TextValue $ws.Range("D15") "18.41"
$ws.Range("E15").Value = "  +2.57%  "

Set-TextValue $ws.Range("D16") "2.991.63"
$ws.Range("E16").Value = "  +1.68%  "

Set-TextValue $ws.Range("D17") "11.19"
$ws.Range("E17").Value = "  +0.40%  "

Set-TextValue $ws.Range("D18") "1.00"
$ws.Range("E18").Value = "  +3.22%  "

Set-TextValue $ws.Range("D19") "51.362.72"
$ws.Range("E19").Value = "  +0.71%  "

Set-TextValue $ws.Range("D20") "3.14"
$ws.Range("E20").Value = "  -0.11%  "

Set-TextValue $ws.Range("D21") "12.56"
$ws.Range("E21").Value = "  +1.67%  "

$dVal = "0.0" + $sub3 + "0961"
Set-TextValue $ws.Range("D22") $dVal
$ws.Range("E22").Value = "  +0.82%  "

Set-TextValue $ws.Range("D23") "70.30"
$ws.Range("E23").Value = "  +2.59%  "

Set-TextValue $ws.Range("D24") "267.09"
$ws.Range("E24").Value = "  +1.13%  "

Set-TextValue $ws.Range("D25") "3.22"
$ws.Range("E25").Value = "  +3.09%  "

Set-TextValue $ws.Range("D26") "7.88"
$ws.Range("E26").Value = "  -2.59%  "

Set-TextValue $ws.Range("D27") "7.53"
$ws.Range("E27").Value = "  +1.35%  "

Set-TextValue $ws.Range("D28") "1.00"
$ws.Range("E28").Value = "  -0.04%  "

Set-TextValue $ws.Range("D29") "26.07"
$ws.Range("E29").Value = "  +2.09%  "

$ws.Range("E30").Value = "  +1.19%  "

$ws.Range("E31").Value = "  -1.76%  "

Set-TextValue $ws.Range("D32") "10.35"
$ws.Range("E32").Value = "  +4.34%  "

Set-TextValue $ws.Range("D33") "34.79"
$ws.Range("E33").Value = "  +5.14%  "

Set-TextValue $ws.Range("D34") "51.47"
$ws.Range("E34").Value = "  +1.68%  "

$ws.Range("E35").Value = "  +1.36%  "

Set-TextValue $ws.Range("D36") "0.0439"
$ws.Range("E36").Value = "  -0.36%  "

$ws.Range("E37").Value = "  -0.07%  "

Set-TextValue $ws.Range("D38") "3.26"
$ws.Range("E38").Value = "  +4.07%  "

$ws.Range("E39").Value = "  +1.56%  "

Set-TextValue $ws.Range("D40") "16.69"
$ws.Range("E40").Value = "  +2.26%  "

$ws.Range("E41").Value = "  +3.57%  "

$ws.Range("E42").Value = "  +3.04%  "

Set-TextValue $ws.Range("D43") "124.39"
$ws.Range("E43").Value = "  +3.96%  "

Set-TextValue $ws.Range("D44") "3.65"
$ws.Range("E44").Value = "  +10.52%  "

Set-TextValue $ws.Range("D45") "21.63"
$ws.Range("E45").Value = "  +2.71%  "

$ws.Range("E46").Value = "  +0.28%  "

Set-TextValue $ws.Range("D47") "2.38"
$ws.Range("E47").Value = "  +4.09%  "

Set-TextValue $ws.Range("D48") "0.271"
$ws.Range("E48").Value = "  -0.77%  "

Set-TextValue $ws.Range("D49") "2.036.59"
$ws.Range("E49").Value = "  +2.57%  "

Set-TextValue $ws.Range("D50") "0.0330"
$ws.Range("E50").Value = "  +1.39%  "

Set-TextValue $ws.Range("D51") "0.529"
$ws.Range("E51").Value = "  +14.86%  "
